# Add a new "Penalty Issued" column (I) to the Dairy Test Threshold report
# template, mirroring the existing header/body-row pattern used for the
# other templated columns (A-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I content -------------------------------------------------
# Row 4 holds the bolded/centered column headers (style copied from H4,
# the "CRY Threshold 3.7" header cell, so the new header matches the
# existing look: bold, centered, wrapped text).
$ws.Range("I4").Value = "Penalty Issued"
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null

# Row 6 is the first templated data row (mirrors {d.Reg[i]...} placeholders
# already present in A6:H6).
$ws.Range("I6").Value = "{d.Reg[i].PenaltyIssued}"

# Rows 3 and 5 are thin spacer rows that, for columns A-H, already carry a
# bottom border via their existing styles; give the new column I cells in
# those rows the same thin-bottom-border treatment.
$ws.Application.CutCopyMode = $false

$border3 = $ws.Range("I3").Borders.Item(9)
$border3.LineStyle = 1
$border3.Weight = 2

$border5 = $ws.Range("I5").Borders.Item(9)
$border5.LineStyle = 1
$border5.Weight = 2

# --- Column sizing ----------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 10.43

# --- Cosmetic selection, matching the saved file's last active cell --------
$ws.Range("G16").Select() | Out-Null

$null = $wb
